# Update the division-fact answers in the worksheet table to the new
# randomly generated set of problems, per commit "Update master to
# output generated at c8c62b6".
#
# The table cells are addressed directly by (row, column) to avoid any
# ambiguity from duplicate text values appearing in both the "old" and
# "new" sets (e.g. "751÷5=150, 1" is both an old value in one cell and
# a new value in another).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$updates = @(
    @{Row=1;  Col=1; Text="911÷5=182, 1"},
    @{Row=1;  Col=2; Text="513÷3=171, 0"},
    @{Row=1;  Col=3; Text="225÷6=37, 3"},
    @{Row=1;  Col=4; Text="356÷7=50, 6"},
    @{Row=1;  Col=5; Text="297÷6=49, 3"},

    @{Row=5;  Col=1; Text="751÷5=150, 1"},
    @{Row=5;  Col=2; Text="146÷4=36, 2"},
    @{Row=5;  Col=3; Text="178÷6=29, 4"},
    @{Row=5;  Col=4; Text="764÷4=191, 0"},
    @{Row=5;  Col=5; Text="427÷3=142, 1"},

    @{Row=9;  Col=1; Text="725÷9=80, 5"},
    @{Row=9;  Col=2; Text="974÷2=487, 0"},
    @{Row=9;  Col=3; Text="725÷4=181, 1"},
    @{Row=9;  Col=4; Text="913÷8=114, 1"},
    @{Row=9;  Col=5; Text="263÷4=65, 3"},

    @{Row=13; Col=1; Text="947÷6=157, 5"},
    @{Row=13; Col=2; Text="974÷6=162, 2"},
    @{Row=13; Col=3; Text="898÷2=449, 0"},
    @{Row=13; Col=4; Text="645÷5=129, 0"},
    @{Row=13; Col=5; Text="563÷3=187, 2"},

    @{Row=17; Col=1; Text="615÷7=87, 6"},
    @{Row=17; Col=2; Text="797÷5=159, 2"},
    @{Row=17; Col=3; Text="131÷2=65, 1"},
    @{Row=17; Col=4; Text="759÷2=379, 1"},
    @{Row=17; Col=5; Text="280÷3=93, 1"}
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.Text
}
